$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.955.29"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.967.22"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9971"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4830"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2948"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06799"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "109.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.970.80"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07743"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.492"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6973"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "293.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "30.948.70"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.681"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007730"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "2.231.18"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9985"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.643"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.935"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.186"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1073"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.443"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.827"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +19.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7747"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02057"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.458"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.136"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8874"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4467"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9982"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.520"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1278"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.396"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "930.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
